$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'C,"1'
$ws.Range("C3").Value = 'C"2'
$ws.Range("C5").Value = "C3`nC3"
$ws.Range("C5").WrapText = $true

$ws.Range("C3").Select()
